$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (date "Förändrad") from 45189 to 45190 for rows 2..138
for ($row = 2; $row -le 138; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
